# Weekly update: insert two new rows (week of 2022-05-25) at the top of the
# Betarraga price history table, pushing the existing rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 492:493 - everything currently at/after row 492
# (down to row 521) shifts down to rows 494:523, dimension grows to R523.
$ws.Rows("492:493").Insert()

# Row 492 - "Primera" grade entry for the new date.
$ws.Cells.Item(492, 1).Value = 9
$ws.Cells.Item(492, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(492, 3).Value = "Metropolitana"
$ws.Cells.Item(492, 4).Value = 44706
$ws.Cells.Item(492, 5).Value = 13
$ws.Cells.Item(492, 6).Value = 100114014
$ws.Cells.Item(492, 7).Value = "Betarraga"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 15000
$ws.Cells.Item(492, 11).Value = 110
$ws.Cells.Item(492, 12).Value = 120
$ws.Cells.Item(492, 13).Value = 115
$ws.Cells.Item(492, 14).Value = "$/unidad"
$ws.Cells.Item(492, 15).Value = "Región Metropolitana"
$ws.Cells.Item(492, 16).Value = 115
$ws.Cells.Item(492, 17).Value = 1
$ws.Cells.Item(492, 18).Value = "Hortaliza"

# Row 493 - "Segunda" grade entry for the same new date.
$ws.Cells.Item(493, 1).Value = 9
$ws.Cells.Item(493, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(493, 3).Value = "Metropolitana"
$ws.Cells.Item(493, 4).Value = 44706
$ws.Cells.Item(493, 5).Value = 13
$ws.Cells.Item(493, 6).Value = 100114014
$ws.Cells.Item(493, 7).Value = "Betarraga"
$ws.Cells.Item(493, 8).Value = "Sin especificar"
$ws.Cells.Item(493, 9).Value = "Segunda"
$ws.Cells.Item(493, 10).Value = 7000
$ws.Cells.Item(493, 11).Value = 100
$ws.Cells.Item(493, 12).Value = 100
$ws.Cells.Item(493, 13).Value = 100
$ws.Cells.Item(493, 14).Value = "$/unidad"
$ws.Cells.Item(493, 15).Value = "Región Metropolitana"
$ws.Cells.Item(493, 16).Value = 100
$ws.Cells.Item(493, 17).Value = 1
$ws.Cells.Item(493, 18).Value = "Hortaliza"
